# Applies the "Updated cryptos list" data refresh (Wed Aug 28 13:50:40 UTC 2024).
#
# Every row's Price (column D) and Volume(1h) (column E) text is refreshed to
# the latest scrape. Two pairs of rows (31/32 and 34/35 and 48/49) were also
# re-ranked, so their Coin name / Link / Price / Volume cells swap places.
#
# Price cells are stored as literal text in the source workbook (e.g.
# "59.572.99", "1.00", "0.0939") even though several of them look like plain
# numbers. Assigning a numeric-looking string straight to `.Value` lets Excel
# "smart type" it into a real number (and mangle the exact text via floating
# point, e.g. "146.76" -> 146.75999999999999). To keep these as text we
# temporarily force the cell to Text format, assign the value, then restore
# the cell's style so no formatting residue is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

$ws.Range('D2').Value = '59.572.99'
$ws.Range('E2').Value = '  -4.24%  '
$ws.Range('D3').Value = '2.486.66'
$ws.Range('E3').Value = '  -4.66%  '
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue 'D5' '540.62'
$ws.Range('E5').Value = '  -2.42%  '
Set-TextValue 'D6' '146.76'
$ws.Range('E6').Value = '  -5.09%  '
Set-TextValue 'D7' '0.996'
$ws.Range('E7').Value = '  -0.44%  '
Set-TextValue 'D8' '0.576'
$ws.Range('E8').Value = '  -2.73%  '
$ws.Range('D9').Value = '2.513.55'
$ws.Range('E9').Value = '  -3.45%  '
Set-TextValue 'D10' '0.101'
$ws.Range('E10').Value = '  -3.58%  '
$ws.Range('E11').Value = '  -1.25%  '
Set-TextValue 'D12' '5.35'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('E13').Value = '  -1.66%  '
$ws.Range('D14').Value = '2.923.73'
$ws.Range('E14').Value = '  -4.84%  '
Set-TextValue 'D15' '24.35'
$ws.Range('E15').Value = '  -5.14%  '
$ws.Range('D16').Value = '59.445.05'
$ws.Range('E16').Value = '  -4.27%  '
$ws.Range('E17').Value = '  -2.82%  '
$ws.Range('D18').Value = '2.509.83'
$ws.Range('E18').Value = '  -4.03%  '
Set-TextValue 'D19' '11.50'
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('E20').Value = '  -3.55%  '
Set-TextValue 'D21' '326.65'
$ws.Range('E21').Value = '  -4.05%  '
Set-TextValue 'D22' '0.998'
$ws.Range('E22').Value = '  -0.11%  '
Set-TextValue 'D23' '5.79'
$ws.Range('E23').Value = '  -4.72%  '
Set-TextValue 'D24' '61.30'
$ws.Range('E24').Value = '  -2.51%  '
$ws.Range('E25').Value = '  -10.32%  '
$ws.Range('E26').Value = '  -2.89%  '
Set-TextValue 'D27' '0.999'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = '2.619.74'
$ws.Range('E28').Value = '  -4.11%  '
Set-TextValue 'D29' '7.85'
Set-TextValue 'D30' '7.12'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D31' '1.30'
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0₃0785'
$ws.Range('E32').Value = '  -5.10%  '
Set-TextValue 'D33' '1.83'
$ws.Range('E33').Value = '  -3.73%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D34' '159.57'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D35' '0.997'
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E36').Value = '  +0.70%  '
Set-TextValue 'D37' '18.72'
$ws.Range('E37').Value = '  -2.75%  '
Set-TextValue 'D38' '4.49'
$ws.Range('E38').Value = '  -4.23%  '
Set-TextValue 'D39' '1.67'
$ws.Range('E39').Value = '  -3.64%  '
Set-TextValue 'D41' '315.38'
$ws.Range('E41').Value = '  -6.66%  '
Set-TextValue 'D42' '36.70'
$ws.Range('E42').Value = '  -2.49%  '
$ws.Range('E43').Value = '  -3.20%  '
$ws.Range('E44').Value = '  -6.19%  '
$ws.Range('E45').Value = '  -0.30%  '
Set-TextValue 'D46' '0.601'
$ws.Range('E46').Value = '  -1.27%  '
Set-TextValue 'D47' '10.79'
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D48' '125.23'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D49' '0.0939'
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('E50').Value = '  -3.54%  '
$ws.Range('E51').Value = '  -3.47%  '
